$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 3622
$ws.Cells.Item(7, 6).Value = 147
$ws.Cells.Item(8, 6).Value = 2345
$ws.Cells.Item(11, 6).Value = 7729
$ws.Cells.Item(13, 6).Value = 629
$ws.Cells.Item(14, 6).Value = 162
$ws.Cells.Item(16, 6).Value = 1004
$ws.Cells.Item(17, 6).Value = 1544
$ws.Cells.Item(18, 6).Value = 2168
$ws.Cells.Item(20, 6).Value = 224
$ws.Cells.Item(21, 6).Value = 284
$ws.Cells.Item(22, 6).Value = 5
$ws.Cells.Item(23, 6).Value = 1117
$ws.Cells.Item(24, 6).Value = 16
$ws.Cells.Item(25, 6).Value = 785
$ws.Cells.Item(26, 6).Value = 69
$ws.Cells.Item(27, 6).Value = 806
$ws.Cells.Item(28, 6).Value = 1346
$ws.Cells.Item(29, 6).Value = 535
$ws.Cells.Item(30, 6).Value = 682
$ws.Cells.Item(33, 6).Value = 37
$ws.Cells.Item(34, 6).Value = 76
$ws.Cells.Item(36, 6).Value = 2528

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 7875
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 23
$ws.Cells.Item(16, 6).Value = 21
$ws.Cells.Item(22, 6).Value = 138
$ws.Cells.Item(31, 6).Value = 93
$ws.Cells.Item(41, 6).Value = 189

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 2441
$ws.Cells.Item(5, 6).Value = 1651
$ws.Cells.Item(7, 6).Value = 702
$ws.Cells.Item(10, 6).Value = 1836
$ws.Cells.Item(16, 6).Value = 2653
$ws.Cells.Item(18, 6).Value = 142
$ws.Cells.Item(19, 6).Value = 595

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 3622
$ws.Cells.Item(4, 6).Value = 2441
$ws.Cells.Item(5, 6).Value = 702
$ws.Cells.Item(6, 6).Value = 1836
$ws.Cells.Item(8, 6).Value = 147
$ws.Cells.Item(12, 6).Value = 7729
$ws.Cells.Item(13, 6).Value = 162
$ws.Cells.Item(14, 6).Value = 142
$ws.Cells.Item(15, 6).Value = 1004
$ws.Cells.Item(16, 6).Value = 1544
$ws.Cells.Item(17, 6).Value = 2168
$ws.Cells.Item(19, 6).Value = 595
$ws.Cells.Item(20, 6).Value = 595
$ws.Cells.Item(22, 6).Value = 23
$ws.Cells.Item(26, 6).Value = 284
$ws.Cells.Item(27, 6).Value = 5
$ws.Cells.Item(28, 6).Value = 1117
$ws.Cells.Item(29, 6).Value = 69
$ws.Cells.Item(30, 6).Value = 806
$ws.Cells.Item(31, 6).Value = 21
$ws.Cells.Item(32, 6).Value = 1346
$ws.Cells.Item(36, 6).Value = 535
$ws.Cells.Item(40, 6).Value = 682
$ws.Cells.Item(47, 6).Value = 2528
